$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header F1 from "quantity" to "categorie_id"
$ws.Range("F1").Value = "categorie_id"

# Update the active selection to H6 (matches the diff's selection change)
$ws.Range("H6").Select()
